$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Stash the two existing hyperlink-cell styles (the plain "Hyperlink" look
# used by J2, and the bordered "Hyperlink" look used by J3:J5) into scratch
# cells so we can re-apply them after the hyperlinks are rebuilt below.
# ---------------------------------------------------------------------------
$ws.Range("J2").Copy() | Out-Null
$ws.Range("ZZ1").PasteSpecial(-4122) | Out-Null
$ws.Range("J3").Copy() | Out-Null
$ws.Range("ZZ2").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# Drop all existing colab_link hyperlinks so they can be rebuilt in order
# with the refreshed notebook URLs (J2:J5) plus three brand-new rows
# (J6:J8). Clear the old display text too, otherwise Excel keeps showing
# the stale URL text even after the hyperlink target is repointed.
# ---------------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete() | Out-Null
$ws.Range("J2:J8").ClearContents() | Out-Null

$colabUrls = @(
  "https://colab.research.google.com/drive/1pY5I1gX0frCy5178bvgAOv27vClxF5b8?usp=sharing",
  "https://colab.research.google.com/drive/1PpCGfVsALwfuq3YuBXFG9zU6KSHuNVUM?usp=sharing",
  "https://colab.research.google.com/drive/1qqzPxSsejvMxF_CcxyvaNYDV3kNAAM6_?usp=sharing",
  "https://colab.research.google.com/drive/13-BIk3bh8Py3Asd0o9T1TrbL0QxezU6P?usp=sharing",
  "https://colab.research.google.com/drive/1aR7xBhCyoLMo_sfkQeRjwyz7FO6e09CK?usp=sharing",
  "https://colab.research.google.com/drive/1Nu7S8G80OL9gpSeOvowtKZjSQ87XbG9V?usp=sharing",
  "https://colab.research.google.com/drive/1mai31zBKjwTT7VLbxvil28x_p9kbQGIN?usp=sharing"
)

for ($i = 0; $i -lt $colabUrls.Length; $i++) {
  $row = 2 + $i
  $cell = $ws.Cells.Item($row, 10)
  $ws.Hyperlinks.Add($cell, $colabUrls[$i]) | Out-Null
}

# Re-apply the original look: row 2 gets the plain style, rows 3-8 get the
# bordered style (matches the pre-existing J3:J5 formatting).
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("J2").PasteSpecial(-4122) | Out-Null

$ws.Range("ZZ2").Copy() | Out-Null
$ws.Range("J3:J8").PasteSpecial(-4122) | Out-Null

$ws.Range("ZZ1").Clear() | Out-Null
$ws.Range("ZZ2").Clear() | Out-Null

# ---------------------------------------------------------------------------
# The three new schedule rows (6, 7, 8) pick up shorter custom row heights
# once their notebook links are added.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).RowHeight = 16
$ws.Rows.Item(7).RowHeight = 14
$ws.Rows.Item(8).RowHeight = 12

# Move the active selection (also clears the stale frozen topLeftCell scroll
# anchor left over from the previous view).
$ws.Range("L7").Select() | Out-Null
